# Updates Price (D) / Volume(1h) (E) columns of the cryptos sheet with the
# latest snapshot values. Percent cells already contain non-numeric text
# (leading/trailing spaces, "%"), so they stay text automatically. Plain
# numeric-looking prices are written with a leading "'" to force Excel to
# keep them as text, matching the workbook's existing inline-string cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.072.24"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "2.243.81"
$ws.Range("E3").Value = "  -1.34%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'315.91"
$ws.Range("D6").Value = "'99.27"
$ws.Range("E6").Value = "  -6.33%  "
$ws.Range("E7").Value = "  -3.16%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  -6.68%  "
$ws.Range("D10").Value = "'36.31"
$ws.Range("E10").Value = "  -6.17%  "
$ws.Range("E11").Value = "  -2.47%  "
$ws.Range("E12").Value = "  -6.76%  "
$ws.Range("E13").Value = "  -2.81%  "
$ws.Range("D14").Value = "2.584.79"
$ws.Range("E14").Value = "  -1.44%  "
$ws.Range("E15").Value = "  -4.36%  "
$ws.Range("D16").Value = "2.250.06"
$ws.Range("E16").Value = "  -1.46%  "
$ws.Range("E17").Value = "  -4.34%  "
$ws.Range("D18").Value = "43.923.66"
$ws.Range("E18").Value = "  -1.04%  "
$ws.Range("E19").Value = "  -6.50%  "
$ws.Range("D20").Value = "0.0₃0980"
$ws.Range("E20").Value = "  -2.61%  "
$ws.Range("D21").Value = "'6.33"
$ws.Range("E21").Value = "  -3.15%  "
$ws.Range("D22").Value = "'65.71"
$ws.Range("E22").Value = "  -1.26%  "
$ws.Range("D23").Value = "'237.93"
$ws.Range("E23").Value = "  -0.64%  "
$ws.Range("E24").Value = "  -7.32%  "
$ws.Range("D25").Value = "'2.03"
$ws.Range("E25").Value = "  -8.20%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "'10.15"
$ws.Range("E27").Value = "  -0.56%  "
$ws.Range("D28").Value = "'2.13"
$ws.Range("E28").Value = "  -4.41%  "
$ws.Range("D29").Value = "'36.45"
$ws.Range("E29").Value = "  -4.88%  "
$ws.Range("D30").Value = "'6.00"
$ws.Range("E30").Value = "  -8.25%  "
$ws.Range("D31").Value = "'20.09"
$ws.Range("E31").Value = "  -2.77%  "
$ws.Range("D32").Value = "'156.31"
$ws.Range("E32").Value = "  -4.74%  "
$ws.Range("D33").Value = "'0.0843"
$ws.Range("E33").Value = "  -5.00%  "
$ws.Range("E34").Value = "  +3.79%  "
$ws.Range("E35").Value = "  -3.46%  "
$ws.Range("D36").Value = "'1.91"
$ws.Range("E36").Value = "  -6.68%  "
$ws.Range("E37").Value = "  -7.57%  "
$ws.Range("E38").Value = "  -3.10%  "
$ws.Range("D39").Value = "'15.43"
$ws.Range("E40").Value = "  -11.27%  "
$ws.Range("D41").Value = "'4.00"
$ws.Range("E41").Value = "  -10.66%  "
$ws.Range("E42").Value = "  -6.23%  "
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").Value = "1.705.89"
$ws.Range("E44").Value = "  -4.36%  "
$ws.Range("D45").Value = "'82.83"
$ws.Range("E45").Value = "  -4.93%  "
$ws.Range("D46").Value = "'0.196"
$ws.Range("E46").Value = "  -6.28%  "
$ws.Range("E47").Value = "  -5.65%  "
$ws.Range("D48").Value = "'102.01"
$ws.Range("E48").Value = "  -2.34%  "
$ws.Range("D49").Value = "'71.38"
$ws.Range("E49").Value = "  -4.82%  "
$ws.Range("D50").Value = "'56.52"
$ws.Range("E50").Value = "  -6.38%  "
$ws.Range("E51").Value = "  -5.93%  "
